$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    394 = 110.53
    395 = 721.79
    396 = 716.18
    397 = 667.05
    398 = 634.12
    399 = 126.58
    400 = 40.99
    401 = 795.05
    402 = 715.63
    403 = 711.3
    404 = 672.28
    405 = 625.4
    406 = 126.91
    407 = 41.07
    408 = 114.56
    409 = 705.21
    410 = 696.63
    411 = 663.53
    412 = 639.64
    413 = 126.89
    414 = 72.86
    415 = 817.05
    416 = 705.08
    417 = 744.02
    418 = 651.68
    419 = 604.9
    420 = 133.04
    421 = 38.25
    422 = 832.59
    423 = 727.77
    424 = 736.04
    425 = 678.65
    426 = 648.68
    427 = 125.79
    428 = 42.5
    429 = 801.36
    430 = 738.12
    431 = 733.96
    432 = 676.95
    433 = 734.35
    434 = 125.84
    435 = 41.45
    436 = 1074.45
    437 = 784.76
    438 = 729.95
    439 = 691.33
    440 = 637.4
    441 = 119.65
    442 = 43.44
    443 = 837.28
    444 = 740.52
    445 = 730.82
    446 = 674.73
    447 = 607.74
    448 = 122.41
    449 = 43.48
    450 = 799.27
    451 = 741.39
    452 = 732.08
    453 = 641.66
    454 = 669.06
    455 = 123.2
    456 = 40.82
    457 = 843.58
    458 = 746.85
    459 = 699.4
    460 = 726.7
    461 = 623.28
    462 = 124.85
    463 = 42.22
    464 = 1004.77
    465 = 739.78
    466 = 727.51
    467 = 726.1
    468 = 623.82
    469 = 121.45
    470 = 41.89
    471 = 1004.06
    472 = 741.36
    473 = 750.2
    474 = 719.51
    475 = 639.16
    476 = 121.41
    477 = 43.55
    478 = 1012.34
    479 = 747.12
    480 = 752.06
    481 = 728.58
    482 = 726.81
    483 = 123.92
    484 = 35.72
    485 = 835.69
    486 = 785.42
    487 = 752.35
    488 = 753.58
    489 = 648.86
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 4).Value = $values[$row]
}

$ws.Range("D394:D489").Select()
